$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MdKbD525"
$ws.Range("B2").Value = 23071985
$ws.Range("C2").Value = "vifsyvt95"
$ws.Range("D2").Value = "QF&8re!9"
$ws.Range("F2").Value = "XWzRTWDw"
$ws.Range("G2").Value = "EWtC"
